# Tygstrup1958.xlsx -- add Hepatic / Systemic / Infusion clearance columns (V, W, X)
# and switch the "P" (galactose elimination rate) column to a 0.00 number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Workbook-level cosmetics (best effort; harmless if the host ignores them)
# ---------------------------------------------------------------------------
$excel.ActiveWindow.TabRatio = 0.141

# ---------------------------------------------------------------------------
# 2) New header (row 3) / unit (row 4) cells for columns V, W, X.
#    Copy the formatting from the existing last column (U) so the new header
#    cells look the same as their neighbours, then fill in the text.
# ---------------------------------------------------------------------------
$ws.Range("U3").Copy()
$ws.Range("V3:X3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("U4").Copy()
$ws.Range("V4:X4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("V3").Value = "Hepatic clearance [ml/min]"
$ws.Range("W3").Value = "Systemic Clearance (CL-CLH) [ml/min]"
$ws.Range("X3").Value = "Clearance via infusion [ml/min]"

$ws.Range("V4").Value = "CLH"
$ws.Range("W4").Value = "CLS"
$ws.Range("X4").Value = "CLI"

# ---------------------------------------------------------------------------
# 3) Data rows 5-38.
#    V = Clearance via infusion as measured upstream ("U" column) normalised
#        per litre of blood   -> U/S*1000
#    W = difference between the galactose infusion rate (J) and V
#    X = Clearance via infusion, computed from the "P" column -> P/S*1000
#
#    Rows whose "U" column only holds the literal "NA" text (no infusion
#    clearance measured) get "NA" for V and W as well; the rest get real
#    formulas.
# ---------------------------------------------------------------------------
$naRows = @(5,6,7,8,9,10,11,12,13,14,30,31,32,36)

for ($r = 5; $r -le 38; $r++) {
    if ($naRows -contains $r) {
        $ws.Cells.Item($r, 22).Value = "NA"
        $ws.Cells.Item($r, 23).Value = "NA"
    } else {
        $ws.Cells.Item($r, 22).Formula = "=U$r/S$r*1000"
        $ws.Cells.Item($r, 23).Formula = "=J$r-V$r"
    }
    $ws.Cells.Item($r, 24).Formula = "=P$r/S$r*1000"
}

# ---------------------------------------------------------------------------
# 4) Number formats: P and the new X column both use "0.00".
# ---------------------------------------------------------------------------
$ws.Range("P5:P38").NumberFormat = "0.00"
$ws.Range("X5:X38").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 5) Selection / scroll position, matching the author's final view.
# ---------------------------------------------------------------------------
$ws.Range("X4").Select()
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1

Write-Host "edit complete"
